$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 243.33333
$ws.Range("I2").Value = 131.25
$ws.Range("J2").Value = 371.42856
$ws.Range("K2").Value = 131.25
$ws.Range("L2").Value = 371.42856
$ws.Range("M2").Value = -18.25
$ws.Range("N2").Value = -597.4285600000001

$ws.Range("H21").Value = 18000
$ws.Range("I21").Value = 18000
$ws.Range("K21").Value = 18000
$ws.Range("M21").Value = -17532

$ws.Range("H23").Value = 18000
$ws.Range("I23").Value = 18000
$ws.Range("K23").Value = 18000
$ws.Range("M23").Value = -17766

$ws.Range("H62").Value = 3482.8333
$ws.Range("I62").Value = 3000
$ws.Range("K62").Value = 3000
$ws.Range("M62").Value = -2376

$ws.Range("H65").Value = 3482.8333
$ws.Range("I65").Value = 3000
$ws.Range("K65").Value = 15000
$ws.Range("M65").Value = -11880

$ws.Range("H103").Value = 1095.6666
$ws.Range("I103").Value = 768.5
$ws.Range("J103").Value = 1750
$ws.Range("K103").Value = 2305.5
$ws.Range("L103").Value = 5250
$ws.Range("M103").Value = -1719.5
$ws.Range("N103").Value = -6422

$ws.Range("H137").Value = 2433
$ws.Range("I137").Value = 1500
$ws.Range("J137").Value = 2772.2727
$ws.Range("K137").Value = 4500
$ws.Range("L137").Value = 8316.8181
$ws.Range("M137").Value = -1950
$ws.Range("N137").Value = -13416.8181

$ws.Range("H138").Value = 2008.5454
$ws.Range("I138").Value = 894.6
$ws.Range("J138").Value = 2133.7078
$ws.Range("K138").Value = 2683.8
$ws.Range("L138").Value = 6401.1234
$ws.Range("M138").Value = 2456.2
$ws.Range("N138").Value = -16681.1234

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1082.6666
$ws.Range("I61").Value = 768.5714
$ws.Range("J61").Value = 2182
$ws.Range("K61").Value = 768.5714
$ws.Range("L61").Value = 2182
$ws.Range("M61").Value = -556.5714
$ws.Range("N61").Value = -2606

$ws.Range("H102").Value = 27794956
$ws.Range("I102").Value = 33353786
$ws.Range("K102").Value = 33353786
$ws.Range("M102").Value = -33352164

$ws.Range("H122").Value = 1452.2
$ws.Range("I122").Value = 1471.2858
$ws.Range("J122").Value = 1407.6666
$ws.Range("K122").Value = 4413.857400000001
$ws.Range("L122").Value = 4222.9998
$ws.Range("M122").Value = -1963.857400000001
$ws.Range("N122").Value = -9122.9998

$ws.Range("H132").Value = 2301.7856
$ws.Range("I132").Value = 2007.5143
$ws.Range("K132").Value = 6022.5429
$ws.Range("M132").Value = -3492.5429

$ws.Range("H136").Value = 1082.6666
$ws.Range("I136").Value = 768.5714
$ws.Range("J136").Value = 2182
$ws.Range("K136").Value = 2305.7142
$ws.Range("L136").Value = 6546
$ws.Range("M136").Value = 244.2857999999997
$ws.Range("N136").Value = -11646

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 20000
$ws.Range("J106").Value = 20000
$ws.Range("L106").Value = 20000
$ws.Range("N106").Value = -22524

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 701.9231
$ws.Range("I105").Value = 649.8889
$ws.Range("J105").Value = 819
$ws.Range("K105").Value = 649.8889
$ws.Range("L105").Value = 819
$ws.Range("M105").Value = 1097.1111
$ws.Range("N105").Value = -4313

$ws.Range("H107").Value = 650.3333
$ws.Range("I107").Value = 298.2
$ws.Range("J107").Value = 760.375
$ws.Range("K107").Value = 298.2
$ws.Range("L107").Value = 760.375
$ws.Range("M107").Value = 1621.8
$ws.Range("N107").Value = -4600.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 3050

$ws.Range("H27").Value = 3050

$ws.Range("H68").Value = 1766.2354
$ws.Range("I68").Value = 899.8333
$ws.Range("J68").Value = 1951.8928
$ws.Range("K68").Value = 2699.4999
$ws.Range("L68").Value = 5855.678400000001
$ws.Range("M68").Value = -1888.4999
$ws.Range("N68").Value = -7477.678400000001

$ws.Range("H71").Value = 1766.2354
$ws.Range("I71").Value = 899.8333
$ws.Range("J71").Value = 1951.8928
$ws.Range("K71").Value = 8098.4997
$ws.Range("L71").Value = 17567.0352
$ws.Range("M71").Value = -4042.4997
$ws.Range("N71").Value = -25679.0352

$ws.Range("H107").Value = 3691.147
$ws.Range("I107").Value = 647.6957
$ws.Range("J107").Value = 10054.728
$ws.Range("K107").Value = 1943.0871
$ws.Range("L107").Value = 30164.184
$ws.Range("M107").Value = -23.08709999999996
$ws.Range("N107").Value = -34004.18399999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 20493.2
$ws.Range("J92").Value = 20493.2
$ws.Range("L92").Value = 20493.2
$ws.Range("N92").Value = -24237.2

$ws.Range("H122").Value = 3100.5715
$ws.Range("I122").Value = 2100.889
$ws.Range("J122").Value = 4900
$ws.Range("K122").Value = 6302.667
$ws.Range("L122").Value = 14700
$ws.Range("M122").Value = -3852.667
$ws.Range("N122").Value = -19600

$ws.Range("H132").Value = 2731.476
$ws.Range("I132").Value = 2297.1538
$ws.Range("K132").Value = 6891.4614
$ws.Range("M132").Value = -4361.4614

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2000.7142
$ws.Range("I7").Value = 1833.3334
$ws.Range("J7").Value = 3005
$ws.Range("K7").Value = 1833.3334
$ws.Range("L7").Value = 3005
$ws.Range("M7").Value = -1721.3334
$ws.Range("N7").Value = -3229

$ws.Range("H40").Value = 3436.25
$ws.Range("I40").Value = 2870
$ws.Range("J40").Value = 4002.5
$ws.Range("K40").Value = 2870
$ws.Range("L40").Value = 4002.5
$ws.Range("M40").Value = -2734
$ws.Range("N40").Value = -4274.5

$ws.Range("H55").Value = 496.44446
$ws.Range("I55").Value = 92.375
$ws.Range("J55").Value = 819.7
$ws.Range("K55").Value = 92.375
$ws.Range("L55").Value = 819.7
$ws.Range("M55").Value = 80.625
$ws.Range("N55").Value = -1165.7

$ws.Range("H68").Value = 1468.5385
$ws.Range("J68").Value = 1816.1666
$ws.Range("L68").Value = 1816.1666
$ws.Range("N68").Value = -3314.1666

$ws.Range("H71").Value = 1468.5385
$ws.Range("J71").Value = 1816.1666
$ws.Range("L71").Value = 9080.833000000001
$ws.Range("N71").Value = -16568.833

$ws.Range("H122").Value = 18892840
$ws.Range("I122").Value = 28337060
$ws.Range("J122").Value = 4400.4
$ws.Range("K122").Value = 85011180
$ws.Range("L122").Value = 13201.2
$ws.Range("M122").Value = -85008730
$ws.Range("N122").Value = -18101.2

$ws.Range("H126").Value = 2000.7142
$ws.Range("I126").Value = 1833.3334
$ws.Range("J126").Value = 3005
$ws.Range("K126").Value = 5500.0002
$ws.Range("L126").Value = 9015
$ws.Range("M126").Value = -3030.0002
$ws.Range("N126").Value = -13955

$ws.Range("H132").Value = 65325.5
$ws.Range("I132").Value = 3240.6
$ws.Range("J132").Value = 93545.91
$ws.Range("K132").Value = 9721.799999999999
$ws.Range("L132").Value = 280637.73
$ws.Range("M132").Value = -7191.799999999999
$ws.Range("N132").Value = -285697.73

$ws.Range("H140").Value = 39929
$ws.Range("J140").Value = 39929
$ws.Range("L140").Value = 39929
$ws.Range("N140").Value = -50289

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 31000
$ws.Range("J46").Value = 31000
$ws.Range("L46").Value = 31000
$ws.Range("N46").Value = -31462

$ws.Range("H122").Value = 23638940
$ws.Range("I122").Value = 23638940
$ws.Range("K122").Value = 70916820
$ws.Range("M122").Value = -70914370

$ws.Range("H126").Value = 48309856
$ws.Range("I126").Value = 52910656
$ws.Range("K126").Value = 158731968
$ws.Range("M126").Value = -158729498

$ws.Range("H132").Value = 4354.48
$ws.Range("I132").Value = 5103.222
$ws.Range("J132").Value = 2429.1428
$ws.Range("K132").Value = 15309.666
$ws.Range("L132").Value = 7287.428400000001
$ws.Range("M132").Value = -12779.666
$ws.Range("N132").Value = -12347.4284

$ws.Range("H134").Value = 31000
$ws.Range("J134").Value = 31000
$ws.Range("L134").Value = 93000
$ws.Range("N134").Value = -98070
